# Magicodes.IE GUID-import test fixture update (issue #44):
# add two new "ProductId"-style test columns (P, Q) with GUID sample
# values to 产品导入模板.xlsx, then widen those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (P2, then Q2) ------------------------------------------------
$ws.Range("P2").Value = "ProductIdTest1"
$ws.Range("Q2").Value = "ProductIdTest2"

# Column P body (P3:P8) ----------------------------------------------------
$ws.Range("P3").Value = "C2EE3694-959A-4A87-BC8C-4003F6576352"
$ws.Range("P4").Value = "C2EE3694-959A-4A87-BC8C-4003F6576353"
$ws.Range("P5").Value = "C2EE3694-959A-4A87-BC8C-4003F6576354"
$ws.Range("P6").Value = "C2EE3694-959A-4A87-BC8C-4003F6576355"
$ws.Range("P7").Value = "C2EE3694-959A-4A87-BC8C-4003F6576356"
$ws.Range("P8").Value = "C2EE3694-959A-4A87-BC8C-4003F6576357"

# Column Q body (Q3:Q5, rows 6-8 intentionally left blank) -----------------
$ws.Range("Q3").Value = "C2EE3694-959A-4A87-BC8C-4003F6576357"
$ws.Range("Q4").Value = "C2EE3694-959A-4A87-BC8C-4003F6576358"
$ws.Range("Q5").Value = "C2EE3694-959A-4A87-BC8C-4003F6576359"

# Widen the two new columns. The engine's ColumnWidth setter stores the
# OOXML <col width> as (chars + 5/7); back that padding out so the
# persisted width lands on exactly 40 / 38 characters.
$ws.Columns.Item(16).ColumnWidth = 40 - 5/7
$ws.Columns.Item(17).ColumnWidth = 38 - 5/7

# Leave the selection where the author ended up after data entry.
$ws.Range("Q11").Select()
